$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: was "Taser [stun]" -> now "Taser [stun] (extra)" + reward note in D4 ---
$ws.Range("A4").Value = "Taser [stun] (extra)"
$ws.Range("D4").Value = 'reward from "Zues"'

# --- Row 5: was "Water [glass]" -> now "Water [glass] (extra)" ---
$ws.Range("A5").Value = "Water [glass] (extra)"

# --- Row 6: was "Lighter [flame]" -> now "Lighter [flame] (extra)" ---
$ws.Range("A6").Value = "Lighter [flame] (extra)"

# --- Rows 19, 20, 22, 23: add "Done" marker in column B ---
$ws.Range("B19").Value = "Done"
$ws.Range("B20").Value = "Done"
$ws.Range("B22").Value = "Done"
$ws.Range("B23").Value = "Done"

# --- New rows 25-27: Enemies section ---
# Copy the formatting of an existing section header ("Mechanics:" at A18) onto
# the new section header cell so it picks up the same fill/shading style.
$ws.Range("A18").Copy()
$ws.Range("A25").PasteSpecial(-4122)
$ws.Range("A25").Value = "Enemies:"

$ws.Range("A26").Value = "Basic Contestant"

$ws.Range("A27").Value = "Zues (Extra)"
$ws.Range("D27").Value = "boss"

# --- Column A widened to fit the new longer entries ---
$ws.Columns("A").ColumnWidth = 19.7

# --- Selection moved to E9:E10 ---
$ws.Range("E9:E10").Select() | Out-Null
